# Apply the "Completed 94 pages of Arroj-e-Iqbal by Azka" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DailyProgress")

# --- Fix J2 (Total no. of pages per day for row 2): 47 -> 72 ---
$ws.Range("J2").Value = 72

# --- Row 3 gets a Resource Name in B3 (it was previously blank) ---
$ws.Range("B3").Value = "Azka Tariq"

# --- Row 3 no longer carries its own "Total no. of pages per day" value ---
$ws.Range("J3").Value = ""

# --- New row 4: progress entry for Arooj-e-Iqbal, 94 pages processed ---
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Range("A4").Value = "23/09/2022"
$ws.Range("B4").Value = "Azka Tariq"
$ws.Range("C4").Value = "Arooj-e-Iqbal"
$ws.Range("D4").Value = 26
$ws.Range("E4").Value = 119
$ws.Range("F4").Value = 94

$ws.Range("G4").Value = 0.85416666666666663
$ws.Range("G4").NumberFormat = "h:mm AM/PM"

$ws.Range("H4").Value = 0.91666666666666663
$ws.Range("H4").NumberFormat = "h:mm AM/PM"

$ws.Range("I4").Value = "1.5 hours"
$ws.Range("J4").Value = 94

# --- Update the view state to match the saved workbook ---
$ws.Range("I14").Select()
$excel.ActiveWindow.ScrollColumn = 3
